$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param([string]$Address, [string]$Text)
    $rng = $ws.Range($Address)
    $rng.NumberFormat = "@"
    $rng.Value = $Text
    $rng.Style = "Normal"
}

# Row 2
Set-TextValue 'D2' '29.954.96'
Set-TextValue 'E2' '  +0.25%  '

# Row 3
Set-TextValue 'D3' '1.882.34'
Set-TextValue 'E3' '  -0.35%  '

# Row 4
Set-TextValue 'D4' '0.9999'

# Row 5
Set-TextValue 'D5' '0.7448'
Set-TextValue 'E5' '  -3.05%  '

# Row 6
Set-TextValue 'D6' '242.98'
Set-TextValue 'E6' '  +0.17%  '

# Row 7
Set-TextValue 'D7' '0.9997'
Set-TextValue 'E7' '  -0.08%  '

# Row 8
Set-TextValue 'D8' '0.3160'
Set-TextValue 'E8' '  +0.90%  '

# Row 9
Set-TextValue 'D9' '0.07221'
Set-TextValue 'E9' '  +1.35%  '

# Row 10
Set-TextValue 'D10' '24.89'
Set-TextValue 'E10' '  -2.77%  '

# Row 11
Set-TextValue 'D11' '0.08339'
Set-TextValue 'E11' '  -2.27%  '

# Row 12
Set-TextValue 'B12' 'Polygon'
Set-TextValue 'C12' 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue 'D12' '0.7552'
Set-TextValue 'E12' '  -1.07%  '

# Row 13
Set-TextValue 'B13' 'Polkadot'
Set-TextValue 'C13' 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue 'D13' '5.419'
Set-TextValue 'E13' '  +1.01%  '

# Row 14
Set-TextValue 'B14' 'WrappedEther'
Set-TextValue 'C14' 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue 'D14' '1.799.89'
Set-TextValue 'E14' '  -5.06%  '

# Row 15
Set-TextValue 'D15' '92.52'
Set-TextValue 'E15' '  -1.16%  '

# Row 16
Set-TextValue 'D16' '6.169'
Set-TextValue 'E16' '  +0.42%  '

# Row 17
Set-TextValue 'D17' '29.996.28'
Set-TextValue 'E17' '  +0.18%  '

# Row 18
Set-TextValue 'D18' '249.52'
Set-TextValue 'E18' '  +2.22%  '

# Row 19
Set-TextValue 'D19' '13.60'
Set-TextValue 'E19' '  -1.08%  '

# Row 20
Set-TextValue 'D20' '0.000007851'
Set-TextValue 'E20' '  +0.37%  '

# Row 21
Set-TextValue 'D21' '2.173.56'
Set-TextValue 'E21' '  +0.82%  '

# Row 22
Set-TextValue 'D22' '1.001'
Set-TextValue 'E22' '  +0.12%  '

# Row 23
Set-TextValue 'D23' '8.013'
Set-TextValue 'E23' '  +0.06%  '

# Row 24
Set-TextValue 'D24' '0.9998'

# Row 25
Set-TextValue 'D25' '0.1563'
Set-TextValue 'E25' '  -4.28%  '

# Row 26
Set-TextValue 'D26' '9.289'
Set-TextValue 'E26' '  -1.01%  '

# Row 27
Set-TextValue 'D27' '165.45'
Set-TextValue 'E27' '  +1.44%  '

# Row 28
Set-TextValue 'D28' '18.71'
Set-TextValue 'E28' '  -0.43%  '

# Row 29
Set-TextValue 'E29' '  +0.29%  '

# Row 30
Set-TextValue 'D30' '1.487'
Set-TextValue 'E30' '  -2.09%  '

# Row 31
Set-TextValue 'D31' '4.614'
Set-TextValue 'E31' '  +2.43%  '

# Row 32
Set-TextValue 'E32' '  +0.07%  '

# Row 33
Set-TextValue 'D33' '4.235'
Set-TextValue 'E33' '  +2.74%  '

# Row 34
Set-TextValue 'D34' '0.05373'
Set-TextValue 'E34' '  -1.38%  '

# Row 35
Set-TextValue 'D35' '1.254'
Set-TextValue 'E35' '  +0.87%  '

# Row 36
Set-TextValue 'D36' '0.7573'
Set-TextValue 'E36' '  +1.53%  '

# Row 37
Set-TextValue 'D37' '0.9925'
Set-TextValue 'E37' '  -0.82%  '

# Row 38
Set-TextValue 'D38' '2.703'
Set-TextValue 'E38' '  -0.03%  '

# Row 39
Set-TextValue 'D39' '0.01967'
Set-TextValue 'E39' '  +0.89%  '

# Row 40
Set-TextValue 'E40' '  -0.64%  '

# Row 41
Set-TextValue 'D41' '0.4560'
Set-TextValue 'E41' '  +1.93%  '

# Row 42
Set-TextValue 'D42' '1.108.01'
Set-TextValue 'E42' '  +0.68%  '

# Row 43
Set-TextValue 'D43' '6.058'
Set-TextValue 'E43' '  -0.36%  '

# Row 44
Set-TextValue 'D44' '72.54'
Set-TextValue 'E44' '  -0.80%  '

# Row 45
Set-TextValue 'D45' '0.8704'
Set-TextValue 'E45' '  +1.56%  '

# Row 46
Set-TextValue 'D46' '104.65'
Set-TextValue 'E46' '  +1.56%  '

# Row 47
Set-TextValue 'D47' '0.9995'
Set-TextValue 'E47' '  -0.10%  '

# Row 48
Set-TextValue 'D48' '1.867'
Set-TextValue 'E48' '  -0.22%  '

# Row 49
Set-TextValue 'D49' '7.611'
Set-TextValue 'E49' '  -0.80%  '

# Row 50
Set-TextValue 'D50' '2.104.66'
Set-TextValue 'E50' '  +3.42%  '

# Row 51
Set-TextValue 'D51' '9.552'
Set-TextValue 'E51' '  -1.91%  '
